# Commit: "change metadata sheet to isa template"
#
# The Swate-template metadata worksheet is renamed from
# "SwateTemplateMetadata" to "isa_template". Matching the workbook's
# recorded active-tab change (tabSelected moves onto this sheet /
# workbookView activeTab="1"), the renamed sheet is also activated.

$wb = $excel.ActiveWorkbook

$metaSheet = $wb.Worksheets.Item("SwateTemplateMetadata")
$metaSheet.Name = "isa_template"
$metaSheet.Activate()
